# Weekly update: insert a new daily price record for Zapallo italiano
# (Vega Central Mapocho de Santiago) at row 528 and push the existing
# rows 528:599 down to 529:600.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 528 (shifts 528:599 -> 529:600,
# and extends the used range / dimension to row 600).
$ws.Rows.Item(528).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(528, 1).Value  = 9
$ws.Cells.Item(528, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(528, 3).Value  = "Metropolitana"
$ws.Cells.Item(528, 4).Value  = 45127
$ws.Cells.Item(528, 5).Value  = 13
$ws.Cells.Item(528, 6).Value  = 100112032
$ws.Cells.Item(528, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(528, 8).Value  = "Sin especificar"
$ws.Cells.Item(528, 9).Value  = "Primera"
$ws.Cells.Item(528, 10).Value = 70
$ws.Cells.Item(528, 11).Value = 16000
$ws.Cells.Item(528, 12).Value = 18000
$ws.Cells.Item(528, 13).Value = 17000
$ws.Cells.Item(528, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(528, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(528, 16).Value = 340
$ws.Cells.Item(528, 17).Value = 50
$ws.Cells.Item(528, 18).Value = "Hortaliza"
